$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.272.63'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '2.508.46'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '321.76'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '108.33'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.66%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -0.30%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '39.92'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.40%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.31'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +8.64%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0819'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("E13").Value = '  +0.01%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.18'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("D15").Value = '2.900.10'
$ws.Range("D16").Value = '2.508.39'
$ws.Range("E16").Value = '  +0.66%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.844'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = '48.122.15'
$ws.Range("E18").Value = '  +1.43%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.14'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("E20").Value = '  +2.28%  '
$ws.Range("D21").Value = '0.0₃0944'
$ws.Range("E21").Value = '  +0.26%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.76'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '72.34'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.33%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '278.30'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +12.69%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -0.08%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '25.76'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.79'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.141'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '35.31'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '2.11'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.41%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '49.53'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.80%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '19.68'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.34'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("E36").Value = '  -0.65%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.96'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.29%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '4.66'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("E40").Value = '  -0.12%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '122.47'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.86%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.22'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.17%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '21.57'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").Value = '2.016.80'
$ws.Range("E45").Value = '  +0.99%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.17'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +4.52%  '
$ws.Range("E47").Value = '  +3.19%  '
$ws.Range("E48").Value = '  -2.29%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '9.05'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("E50").Value = '  -0.71%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '80.58'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +3.75%  '
